# Updates to Sample documents
#
# Three WSP package names in the "Solution Assessment Report" sample deck
# are corrected:
#   - Slide 13: "contoso.sharepoint.libraryreceivers.wsp" -> "contoso.sharepoint.eventreceivers.wsp"
#   - Slide 17: "contoso.sharepoint.receivers.wsp"         -> "contoso.sharepoint.eventreceivers.wsp"
#   - Slide 18: "fabricam.locationfinder.wsp"               -> "fabrikam.locationfinder.wsp"

$p = $ppt.ActivePresentation

# --- Slide 13: "Text Placeholder 2" bullet list, 5th paragraph -------------
$slide13 = $p.Slides.Item(13)
$shape13 = $slide13.Shapes.Item(2)
$run13 = $shape13.TextFrame.TextRange.Paragraphs(5, 1).Runs(1, 1)
$run13.Text = "contoso.sharepoint.eventreceivers.wsp"

# --- Slide 17: "Table 6", row 6 / column 1 ----------------------------------
$slide17 = $p.Slides.Item(17)
$table17 = $slide17.Shapes.Item(2).Table
$run17 = $table17.Cell(6, 1).Shape.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$run17.Text = "contoso.sharepoint.eventreceivers.wsp"

# --- Slide 18: "Table 4", row 4 / column 1 ----------------------------------
$slide18 = $p.Slides.Item(18)
$table18 = $slide18.Shapes.Item(2).Table
$run18 = $table18.Cell(4, 1).Shape.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$run18.Text = "fabrikam.locationfinder.wsp"
